# Add a new row (34) to Sheet1 documenting the GFG "sort linked list of 0s,1s,2s" question.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A34").Value = "GFG"
$ws.Range("A34").HorizontalAlignment = -4108
$ws.Range("A34").VerticalAlignment = -4160

$ws.Range("B34").Value = "Given a linked list of 0s, 1s and 2s, sort it."

$ws.Range("C34").Value = "Java"

# Update selection to mirror the authored state (active cell C34)
$ws.Range("C34").Select()
